# Update due to suppressing oids:
#  - "Experimental" (row 7) on the Metadata sheet was blank; it now holds the
#    literal text "false" (not the boolean FALSE).
#  - "Date" (row 8) on the Metadata sheet is refreshed to the new export
#    timestamp.
#
# Typing the bare word false/FALSE into a cell's .Value is auto-coerced to a
# real Boolean by Excel, which is not what the source workbook stores (it
# keeps "false" as plain text). So the literal text is produced with a
# formula that evaluates to the string "false", then copied/pasted as a
# value (paste-values preserves the text type instead of re-parsing it).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$helper = $ws.Range("D1")
$helper.Formula = "=""fal""&""se"""
$helper.Copy()
$ws.Range("B7").PasteSpecial(-4163)
$helper.ClearContents()

$ws.Range("B8").Value = "2023-10-09T22:41:16+02:00"
